$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 1 closure entry (row 6) ---
$ws.Range("A6").Value = (Get-Date -Year 2021 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B6").Value = "Festivel"
$ws.Range("C6").Value = "Ice Scrum"
$ws.Range("D6").Value = "Clôture du Sprint 1 "

# --- Database link entry (row 7) ---
$ws.Range("A7").Value = (Get-Date -Year 2021 -Month 5 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B7").Value = "Festivel"
$ws.Range("C7").Value = "Code"
$ws.Range("D7").Value = "Liasion à la base de donnée"

# --- move the active selection cursor, matching the author's last position ---
$null = $ws.Range("D8").Select()
